$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1340
$ws1.Range("F8").Value = 11592
$ws1.Range("F15").Value = 1090
$ws1.Range("F18").Value = 4669
$ws1.Range("F28").Value = 17

# Sheet "全部类型" (all types list)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1340
$ws4.Range("F8").Value = 11592
$ws4.Range("F16").Value = 1090
$ws4.Range("F19").Value = 4669
$ws4.Range("F29").Value = 17
